# Insert a new data row at row 552 (pushing existing rows 552:631 down to 553:632)
# and populate the new row with the latest weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A552").EntireRow.Insert()

$ws.Range("A552").Value2 = 10
$ws.Range("B552").Value2 = "Vega Modelo de Temuco"
$ws.Range("C552").Value2 = "La Araucanía"
$ws.Range("D552").Value2 = 45131
$ws.Range("D552").NumberFormat = $ws.Range("D553").NumberFormat
$ws.Range("E552").Value2 = 9
$ws.Range("F552").Value2 = 100114014
$ws.Range("G552").Value2 = "Betarraga"
$ws.Range("H552").Value2 = "Sin especificar"
$ws.Range("I552").Value2 = "Primera"
$ws.Range("J552").Value2 = 110
$ws.Range("K552").Value2 = 8000
$ws.Range("L552").Value2 = 8000
$ws.Range("M552").Value2 = 8000
$ws.Range("N552").Value2 = '$/saco 25 kilos'
$ws.Range("O552").Value2 = "Provincia de Cautín"
$ws.Range("P552").Value2 = 320
$ws.Range("Q552").Value2 = 25
$ws.Range("R552").Value2 = "Hortaliza"
